# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers in AD1, AE1, AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) so the new
# header cells match the formatting (bold, bordered, centered) of the rest
# of the header row.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in the season record values for each data row (2-41)
for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 30).Value = 96  # AD: Wins
    $ws.Cells.Item($row, 31).Value = 66  # AE: Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF: Ties
}
